$d = $word.ActiveDocument

# --- Step 1: merge the "Versi" + "on" runs into a single "Version" run ---
# Setting Range.Text to the text it already logically represents is a no-op
# for the underlying run structure, so first force a real text change, then
# set it back to the desired final text; this collapses the two runs into one.
$r1 = $d.Range(0, 7)
$r1.Text = "VersionX"
$r1b = $d.Range(0, 8)
$r1b.Text = "Version"

# --- Step 2: change the " 2" run to " 1." (this run sits right before the
# _GoBack bookmark, so editing just this run keeps the bookmark intact) ---
$r3 = $d.Range(7, 9)
$r3.Text = " 1."

# --- Step 3: remove the now-redundant trailing "." run (was after the
# bookmark) without touching the bookmark itself ---
$r4 = $d.Range(10, 11)
$r4.Delete()
